$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New log entries (rows 66-68) -----------------------------------------
# Shared strings are appended to xl/sharedStrings.xml in the order cells are
# first written, so the cells below are populated in the same order the
# original author typed them (per the authoritative diff) to reproduce the
# exact shared-string table ordering.

# B66 "4:42PM"
$ws.Range("B66").Value = "4:42PM"

# G67 "Build simpler invoice gui => for time's sake"
$ws.Range("G67").Value = "Build simpler invoice gui => for time's sake"

# H67 notes
$ws.Range("H67").Value = "Only allow a single invoice to be selected from all unpaid invoices, instead of sending an invoice w/ ability to batch orders together"

# F66 "Create GUI"
$ws.Range("F66").Value = "Create GUI"

# C66 "4:56PM"
$ws.Range("C66").Value = "4:56PM"

# G66 "Create GUI in Qt Designer"
$ws.Range("G66").Value = "Create GUI in Qt Designer"

# H66 notes
$ws.Range("H66").Value = "now need to import to program code"

# B67 "4:57PM"
$ws.Range("B67").Value = "4:57PM"

# C67 "5:09PM" (B68 reuses the same string)
$ws.Range("C67").Value = "5:09PM"
$ws.Range("B68").Value = "5:09PM"

# G68 notes
$ws.Range("G68").Value = "Debug GUI for newInvoiceCandS not showing correctly"

# Cells reusing already-existing shared strings / plain numbers
$ws.Range("F67").Value = "Code"
$ws.Range("F68").Value = "Debug"
$ws.Range("E66").Value = 14

# Date for row 66
$ws.Range("A66").Value = 43810
$ws.Range("A66").NumberFormat = "m/d/yy"

# A67 carries the date's style but stays empty
$ws.Range("A67").NumberFormat = "m/d/yy"

# Time-style (h:mm) formatting applied to C66, C67 and B68, matching style "4"
$ws.Range("C66").NumberFormat = "h:mm"
$ws.Range("C67").NumberFormat = "h:mm"
$ws.Range("B68").NumberFormat = "h:mm"

# --- View cosmetic updates --------------------------------------------------
# Scroll the viewport towards the new rows, then land the final selection on
# G68 (the last edited cell), matching the saved sheetView/selection state.
$ws.Range("A28").Select()
$ws.Range("G68").Select()
